$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.455198407173157
$ws.Range("B1").Value = 1.777684688568115
$ws.Range("C1").Value = 2.447533369064331
$ws.Range("D1").Value = 4.40050220489502
$ws.Range("E1").Value = 2.47014594078064
